# Auto-generated edit script applying the Ravana_Profits.xlsx market-data refresh
# described by the commit diff. Values are plain numbers (no formulas in source).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # row 8: set "H8"=56.5,"I8"=56.5,"K8"=169.5,"M8"=-30.5
$ws.Range("H8").Value = 56.5
$ws.Range("I8").Value = 56.5
$ws.Range("K8").Value = 169.5
$ws.Range("M8").Value = -30.5

  # row 94: set "H94"=13005,"I94"=13005,"K94"=13005,"M94"=-12554
$ws.Range("H94").Value = 13005
$ws.Range("I94").Value = 13005
$ws.Range("K94").Value = 13005
$ws.Range("M94").Value = -12554

  # row 137: set "H137"=4441.1055,"I137"=2555.4285,"K137"=7666.2855,"M137"=-5116.2855
$ws.Range("H137").Value = 4441.1055
$ws.Range("I137").Value = 2555.4285
$ws.Range("K137").Value = 7666.2855
$ws.Range("M137").Value = -5116.2855

$ws = $wb.Worksheets.Item("ARM")
  # row 2: set "H2"=1203.6666,"I2"=805.5,"K2"=805.5,"M2"=-692.5
$ws.Range("H2").Value = 1203.6666
$ws.Range("I2").Value = 805.5
$ws.Range("K2").Value = 805.5
$ws.Range("M2").Value = -692.5

  # row 116: set "H116"=1203.6666,"I116"=805.5,"K116"=805.5,"M116"=1488.5
$ws.Range("H116").Value = 1203.6666
$ws.Range("I116").Value = 805.5
$ws.Range("K116").Value = 805.5
$ws.Range("M116").Value = 1488.5

$ws = $wb.Worksheets.Item("BSM")
  # row 3: set "H3"=1203.6666,"I3"=805.5,"K3"=805.5,"M3"=-691.5
$ws.Range("H3").Value = 1203.6666
$ws.Range("I3").Value = 805.5
$ws.Range("K3").Value = 805.5
$ws.Range("M3").Value = -691.5

  # row 86: set "H86"=2850,"I86"=0,"J86"=2850,"K86"=0,"M86"=2850,"N86"=-5096; clear L86
$ws.Range("L86").ClearContents()
$ws.Range("H86").Value = 2850
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2850
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = 2850
$ws.Range("N86").Value = -5096

  # row 89: set "H89"=2850,"I89"=0,"J89"=2850,"K89"=0,"M89"=14250,"N89"=-25482; clear L89
$ws.Range("L89").ClearContents()
$ws.Range("H89").Value = 2850
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2850
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = 14250
$ws.Range("N89").Value = -25482

  # row 134: set "H134"=3805.7693,"I134"=3805.7693,"J134"=0,"K134"=11417.3079,"L134"=0,"N134"=-8882.3079; clear M134
$ws.Range("M134").ClearContents()
$ws.Range("H134").Value = 3805.7693
$ws.Range("I134").Value = 3805.7693
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11417.3079
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -8882.3079

$ws = $wb.Worksheets.Item("CRP")
  # row 58: set "H58"=3249.3333,"I58"=3299.4,"J58"=2999,"K58"=3299.4,"L58"=2999,"M58"=-3096.4,"N58"=-3405
$ws.Range("H58").Value = 3249.3333
$ws.Range("I58").Value = 3299.4
$ws.Range("J58").Value = 2999
$ws.Range("K58").Value = 3299.4
$ws.Range("L58").Value = 2999
$ws.Range("M58").Value = -3096.4
$ws.Range("N58").Value = -3405

  # row 86: set "H86"=21519.312,"I86"=5998.25,"J86"=37040.375,"K86"=5998.25,"L86"=37040.375,"M86"=-4875.25,"N86"=-39286.375
$ws.Range("H86").Value = 21519.312
$ws.Range("I86").Value = 5998.25
$ws.Range("J86").Value = 37040.375
$ws.Range("K86").Value = 5998.25
$ws.Range("L86").Value = 37040.375
$ws.Range("M86").Value = -4875.25
$ws.Range("N86").Value = -39286.375

  # row 89: set "H89"=21519.312,"I89"=5998.25,"J89"=37040.375,"K89"=29991.25,"L89"=185201.875,"M89"=-24375.25,"N89"=-196433.875
$ws.Range("H89").Value = 21519.312
$ws.Range("I89").Value = 5998.25
$ws.Range("J89").Value = 37040.375
$ws.Range("K89").Value = 29991.25
$ws.Range("L89").Value = 185201.875
$ws.Range("M89").Value = -24375.25
$ws.Range("N89").Value = -196433.875

  # row 99: set "H99"=5000,"I99"=5000,"K99"=5000,"M99"=-3502
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -3502

  # row 122: set "H122"=1878.8462,"I122"=1947.8182,"K122"=5843.4546,"M122"=-3393.4546
$ws.Range("H122").Value = 1878.8462
$ws.Range("I122").Value = 1947.8182
$ws.Range("K122").Value = 5843.4546
$ws.Range("M122").Value = -3393.4546

  # row 126: set "H126"=5000,"I126"=5000,"K126"=15000,"M126"=-12530
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

  # row 132: set "H132"=4026.25,"I132"=3727.5,"J132"=4623.75,"K132"=11182.5,"L132"=13871.25,"M132"=-8652.5,"N132"=-18931.25
$ws.Range("H132").Value = 4026.25
$ws.Range("I132").Value = 3727.5
$ws.Range("J132").Value = 4623.75
$ws.Range("K132").Value = 11182.5
$ws.Range("L132").Value = 13871.25
$ws.Range("M132").Value = -8652.5
$ws.Range("N132").Value = -18931.25

  # row 134: set "H134"=5999,"I134"=5999,"K134"=17997,"M134"=-15462
$ws.Range("H134").Value = 5999
$ws.Range("I134").Value = 5999
$ws.Range("K134").Value = 17997
$ws.Range("M134").Value = -15462

  # row 136: set "H136"=3249.3333,"I136"=3299.4,"J136"=2999,"K136"=9898.200000000001,"L136"=8997,"M136"=-7348.200000000001,"N136"=-14097
$ws.Range("H136").Value = 3249.3333
$ws.Range("I136").Value = 3299.4
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 9898.200000000001
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -7348.200000000001
$ws.Range("N136").Value = -14097

$ws = $wb.Worksheets.Item("CUL")
  # row 4: set "H4"=125026770,"I4"=87401230,"K4"=262203690,"M4"=-262203578
$ws.Range("H4").Value = 125026770
$ws.Range("I4").Value = 87401230
$ws.Range("K4").Value = 262203690
$ws.Range("M4").Value = -262203578

  # row 68: set "H68"=1837.25,"I68"=2074.5,"J68"=1600,"K68"=6223.5,"L68"=4800,"M68"=-5412.5,"N68"=-6422
$ws.Range("H68").Value = 1837.25
$ws.Range("I68").Value = 2074.5
$ws.Range("J68").Value = 1600
$ws.Range("K68").Value = 6223.5
$ws.Range("L68").Value = 4800
$ws.Range("M68").Value = -5412.5
$ws.Range("N68").Value = -6422

  # row 69: set "H69"=0,"J69"=0,"N69"=0; clear L69
$ws.Range("L69").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("N69").Value = 0

  # row 71: set "H71"=1837.25,"I71"=2074.5,"J71"=1600,"K71"=18670.5,"L71"=14400,"M71"=-14614.5,"N71"=-22512
$ws.Range("H71").Value = 1837.25
$ws.Range("I71").Value = 2074.5
$ws.Range("J71").Value = 1600
$ws.Range("K71").Value = 18670.5
$ws.Range("L71").Value = 14400
$ws.Range("M71").Value = -14614.5
$ws.Range("N71").Value = -22512

  # row 72: set "H72"=0,"J72"=0,"N72"=0; clear L72
$ws.Range("L72").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("N72").Value = 0

  # row 80: set "H80"=0,"J80"=0,"N80"=0; clear L80
$ws.Range("L80").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").Value = 0

  # row 83: set "H83"=0,"J83"=0,"N83"=0; clear L83
$ws.Range("L83").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("N83").Value = 0

  # row 132: set "H132"=4390.7144,"I132"=3934.75,"K132"=35412.75,"M132"=-32882.75
$ws.Range("H132").Value = 4390.7144
$ws.Range("I132").Value = 3934.75
$ws.Range("K132").Value = 35412.75
$ws.Range("M132").Value = -32882.75

  # row 133: set "H133"=0,"I133"=0,"K133"=0; clear M133
$ws.Range("M133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0

  # row 138: set "H138"=3000,"I138"=3000,"J138"=0,"K138"=9000,"M138"=-3860,"N138"=0; clear L138
$ws.Range("L138").ClearContents()
$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 3000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9000
$ws.Range("M138").Value = -3860
$ws.Range("N138").Value = 0

  # row 140: set "H140"=3341.75,"I140"=3341.75,"K140"=10025.25,"M140"=-4845.25
$ws.Range("H140").Value = 3341.75
$ws.Range("I140").Value = 3341.75
$ws.Range("K140").Value = 10025.25
$ws.Range("M140").Value = -4845.25

$ws = $wb.Worksheets.Item("GSM")
  # row 15: set "H15"=30118,"J15"=30118,"L15"=30118,"N15"=-30694
$ws.Range("H15").Value = 30118
$ws.Range("J15").Value = 30118
$ws.Range("L15").Value = 30118
$ws.Range("N15").Value = -30694

  # row 81: set "H81"=30118,"J81"=30118,"L81"=30118,"N81"=-32114
$ws.Range("H81").Value = 30118
$ws.Range("J81").Value = 30118
$ws.Range("L81").Value = 30118
$ws.Range("N81").Value = -32114

  # row 84: set "H84"=30118,"J84"=30118,"L84"=90354,"N84"=-100338
$ws.Range("H84").Value = 30118
$ws.Range("J84").Value = 30118
$ws.Range("L84").Value = 90354
$ws.Range("N84").Value = -100338

  # row 122: set "H122"=800,"I122"=800,"K122"=2400,"M122"=50
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50

  # row 126: set "H126"=1499.125,"I126"=1499,"K126"=4497,"M126"=-2027
$ws.Range("H126").Value = 1499.125
$ws.Range("I126").Value = 1499
$ws.Range("K126").Value = 4497
$ws.Range("M126").Value = -2027

  # row 132: set "H132"=3000,"I132"=3000,"J132"=3000,"K132"=9000,"L132"=9000,"M132"=-6470,"N132"=-14060
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
  # row 40: set "H40"=5329.6665,"I40"=5244.75,"K40"=5244.75,"M40"=-5108.75
$ws.Range("H40").Value = 5329.6665
$ws.Range("I40").Value = 5244.75
$ws.Range("K40").Value = 5244.75
$ws.Range("M40").Value = -5108.75

  # row 122: set "H122"=2900.75,"I122"=2710.182,"K122"=8130.545999999999,"M122"=-5680.545999999999
$ws.Range("H122").Value = 2900.75
$ws.Range("I122").Value = 2710.182
$ws.Range("K122").Value = 8130.545999999999
$ws.Range("M122").Value = -5680.545999999999

  # row 132: set "H132"=2636.5334,"I132"=1792.1666,"K132"=5376.4998,"M132"=-2846.4998
$ws.Range("H132").Value = 2636.5334
$ws.Range("I132").Value = 1792.1666
$ws.Range("K132").Value = 5376.4998
$ws.Range("M132").Value = -2846.4998

  # row 136: set "H136"=2805.5908,"I136"=2595.4443,"K136"=7786.3329,"M136"=-5236.3329
$ws.Range("H136").Value = 2805.5908
$ws.Range("I136").Value = 2595.4443
$ws.Range("K136").Value = 7786.3329
$ws.Range("M136").Value = -5236.3329

$ws = $wb.Worksheets.Item("WVR")
  # row 8: set "H8"=0,"J8"=0,"N8"=0; clear L8
$ws.Range("L8").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").Value = 0

  # row 87: set "H87"=89999,"J87"=89999,"L87"=89999,"N87"=-92495
$ws.Range("H87").Value = 89999
$ws.Range("J87").Value = 89999
$ws.Range("L87").Value = 89999
$ws.Range("N87").Value = -92495

  # row 90: set "H90"=89999,"J90"=89999,"L90"=269997,"N90"=-282477
$ws.Range("H90").Value = 89999
$ws.Range("J90").Value = 89999
$ws.Range("L90").Value = 269997
$ws.Range("N90").Value = -282477

  # row 122: set "H122"=2244,"I122"=2244,"J122"=0,"K122"=6732,"L122"=0,"N122"=-4282; clear M122
$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 2244
$ws.Range("I122").Value = 2244
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6732
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -4282

  # row 126: set "H126"=2524.6667,"I126"=2556.111,"J126"=2430.3333,"K126"=7668.333,"L126"=7290.999899999999,"M126"=-5198.333,"N126"=-12230.9999
$ws.Range("H126").Value = 2524.6667
$ws.Range("I126").Value = 2556.111
$ws.Range("J126").Value = 2430.3333
$ws.Range("K126").Value = 7668.333
$ws.Range("L126").Value = 7290.999899999999
$ws.Range("M126").Value = -5198.333
$ws.Range("N126").Value = -12230.9999

  # row 136: set "H136"=1163.56,"I136"=984.1905,"K136"=2952.5715,"M136"=-402.5715
$ws.Range("H136").Value = 1163.56
$ws.Range("I136").Value = 984.1905
$ws.Range("K136").Value = 2952.5715
$ws.Range("M136").Value = -402.5715
